$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -1.00176269401673
$ws.Range("C2").Value = 3.50752444718883

$ws.Range("B3").Value = -1.01615465315107
$ws.Range("C3").Value = -1.52295354373166

$ws.Range("B4").Value = -0.40311910543042
$ws.Range("C4").Value = 5.22577178225862

$ws.Range("B5").Value = -0.408410218063001
$ws.Range("C5").Value = 2.41793441867662

$ws.Range("B6").Value = 1.51791677118226
$ws.Range("C6").Value = -1.89206233685472

$ws.Range("B7").Value = -0.545597048759035
$ws.Range("C7").Value = 2.08165152525859

$ws.Range("B8").Value = -0.488182793699075
$ws.Range("C8").Value = 0.26306172346934

$ws.Range("B9").Value = -0.0257851821469864
$ws.Range("C9").Value = 3.6143848167625

$ws.Range("B10").Value = -0.214325197035093
$ws.Range("C10").Value = 3.73719427038262

$ws.Range("B11").Value = -0.513154107203089
$ws.Range("C11").Value = 0.562965673323551
